# Updates cryptos list data (prices / 1h volume %) as scraped on
# Sat Sep  7 23:59:40 UTC 2024 with GitHub Actions.
#
# Helper: write a value into a cell as TEXT (never let Excel reinterpret
# a numeric-looking string like "4.67" as a Number), while keeping the
# cell's style/format exactly as it was (no visible NumberFormat/Style
# change) by resetting the style back to Normal right after.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "54.109.76"
$ws.Range("E2").Value = "  +0.79%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.258.51"
$ws.Range("E3").Value = "  +1.54%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "493.38"
$ws.Range("E5").Value = "  +1.26%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "127.69"
$ws.Range("E6").Value = "  +2.20%  "

# Row 7 - USDC
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.525"
$ws.Range("E8").Value = "  +0.75%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.0953"
$ws.Range("E9").Value = "  +3.00%  "

# Row 10 - TRON
Set-TextValue $ws.Range("D10") "0.152"
$ws.Range("E10").Value = "  +2.51%  "

# Row 11 - Cardano
Set-TextValue $ws.Range("D11") "0.325"
$ws.Range("E11").Value = "  +3.33%  "

# Row 12 - Toncoin
Set-TextValue $ws.Range("D12") "4.67"
$ws.Range("E12").Value = "  +0.34%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "2.661.37"
$ws.Range("E13").Value = "  +1.81%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "21.79"
$ws.Range("E14").Value = "  +3.05%  "

# Row 15 - WrappedBTC
Set-TextValue $ws.Range("D15") "54.010.40"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16 - ShibaInu
Set-TextValue $ws.Range("D16") "0.0000129"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.297.34"
$ws.Range("E17").Value = "  +4.04%  "

# Row 18 - Chainlink
Set-TextValue $ws.Range("D18") "10.02"
$ws.Range("E18").Value = "  +4.73%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.95%  "

# Row 20 - was Uniswap, now BitcoinCash
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "300.56"
$ws.Range("E20").Value = "  +1.88%  "

# Row 21 - was BitcoinCash, now Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D21") "6.42"
$ws.Range("E21").Value = "  +4.07%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.06%  "

# Row 23 - LEO
$ws.Range("E23").Value = "  -1.47%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "62.08"
$ws.Range("E24").Value = "  -1.27%  "

# Row 25 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D25") "0.998"
$ws.Range("E25").Value = "  +0.01%  "

# Row 26 - Polygon
$ws.Range("E26").Value = "  +0.97%  "

# Row 27 - WrappedeETH
Set-TextValue $ws.Range("D27") "2.363.99"
$ws.Range("E27").Value = "  +1.63%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +1.72%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "7.06"
$ws.Range("E29").Value = "  +0.55%  "

# Row 30 - Monero
Set-TextValue $ws.Range("D30") "167.53"
$ws.Range("E30").Value = "  +1.48%  "

# Row 31 - PancakeSwap
Set-TextValue $ws.Range("D31") "1.59"
$ws.Range("E31").Value = "  +0.77%  "

# Row 32 - Aptos
Set-TextValue $ws.Range("D32") "5.86"
$ws.Range("E32").Value = "  +2.71%  "

# Row 33 - PEPE
Set-TextValue $ws.Range("D33") "0.0₃0677"
$ws.Range("E33").Value = "  +1.82%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  -0.01%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.02%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  -0.31%  "

# Row 37 - EthereumClassic
Set-TextValue $ws.Range("D37") "17.66"
$ws.Range("E37").Value = "  +2.12%  "

# Row 38 - SuiNetwork
Set-TextValue $ws.Range("D38") "0.890"
$ws.Range("E38").Value = "  +6.26%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +2.56%  "

# Row 40 - NEARProtocol
Set-TextValue $ws.Range("D40") "3.67"
$ws.Range("E40").Value = "  +3.38%  "

# Row 41 - OKB
Set-TextValue $ws.Range("D41") "35.74"
$ws.Range("E41").Value = "  -0.45%  "

# Row 42 - Stacks
Set-TextValue $ws.Range("D42") "1.39"
$ws.Range("E42").Value = "  +1.75%  "

# Row 43 - PolygonEcosystemToken
$ws.Range("E43").Value = "  +0.82%  "

# Row 44 - Filecoin
Set-TextValue $ws.Range("D44") "3.35"
$ws.Range("E44").Value = "  +2.35%  "

# Row 45 - was Aave, now RenderToken
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") "4.92"
$ws.Range("E45").Value = "  +2.87%  "

# Row 46 - was RenderToken, now Aave
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "125.12"
$ws.Range("E46").Value = "  -1.32%  "

# Row 47 - Stellar
Set-TextValue $ws.Range("D47") "0.0886"
$ws.Range("E47").Value = "  +0.59%  "

# Row 48 - Mantle
Set-TextValue $ws.Range("D48") "0.542"
$ws.Range("E48").Value = "  +1.16%  "

# Row 49 - Hedera
Set-TextValue $ws.Range("D49") "0.0485"
$ws.Range("E49").Value = "  +3.00%  "

# Row 50 - Bittensor
Set-TextValue $ws.Range("D50") "235.10"
$ws.Range("E50").Value = "  +1.55%  "

# Row 51 - VeChain
Set-TextValue $ws.Range("D51") "0.0202"
$ws.Range("E51").Value = "  +1.14%  "
